# Updates the "cryptos" price/volume table to the latest scrape values
# (GitHub Actions refresh). Columns: A=rank, B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings such as "0.9987" or "241.80" parse as plain numbers, and Excel
# would silently convert them away from text. Force the cell to text format,
# assign the value, then drop back to the default "Normal" style so no stray
# number-format style lingers on the cell.
function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '29.503.86'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '1.851.40'
$ws.Range("E3").Value = '  -0.53%  '
Set-TextValue 'D4' '0.9987'
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue 'D5' '241.80'
Set-TextValue 'D6' '0.6298'
$ws.Range("E6").Value = '  -2.58%  '
Set-TextValue 'D7' '0.9996'
$ws.Range("E7").Value = '  -0.02%  '
Set-TextValue 'D8' '0.07536'
$ws.Range("E8").Value = '  -0.08%  '
Set-TextValue 'D9' '0.2980'
$ws.Range("E9").Value = '  -0.39%  '
Set-TextValue 'D10' '24.44'
Set-TextValue 'D11' '0.07723'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").Value = '1.840.02'
$ws.Range("E12").Value = '  -1.30%  '
Set-TextValue 'D13' '0.6927'
$ws.Range("E13").Value = '  -0.12%  '
Set-TextValue 'D14' '5.007'
$ws.Range("E14").Value = '  -1.04%  '
Set-TextValue 'D15' '83.69'
$ws.Range("E15").Value = '  -0.46%  '
Set-TextValue 'D16' '0.000009800'
$ws.Range("E16").Value = '  -0.37%  '
$ws.Range("D17").Value = '2.145.96'
$ws.Range("E17").Value = '  +1.21%  '
Set-TextValue 'D18' '6.240'
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("D19").Value = '29.523.19'
$ws.Range("E19").Value = '  -1.07%  '
Set-TextValue 'D20' '233.43'
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("E21").Value = '  -1.46%  '
Set-TextValue 'D22' '0.9992'
$ws.Range("E22").Value = '  -0.02%  '
Set-TextValue 'D23' '7.657'
$ws.Range("E23").Value = '  +0.26%  '
Set-TextValue 'D24' '0.9997'
$ws.Range("E24").Value = '  -0.01%  '
Set-TextValue 'D25' '154.62'
$ws.Range("E25").Value = '  -2.40%  '
Set-TextValue 'D26' '0.1391'
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("E28").Value = '  -1.32%  '
Set-TextValue 'D29' '1.476'
$ws.Range("E29").Value = '  -1.19%  '
Set-TextValue 'D30' '0.05924'
$ws.Range("E30").Value = '  -4.58%  '
Set-TextValue 'D31' '1.251'
$ws.Range("E31").Value = '  -3.19%  '
Set-TextValue 'D32' '4.107'
$ws.Range("E32").Value = '  -1.36%  '
Set-TextValue 'D33' '4.034'
$ws.Range("E33").Value = '  -1.65%  '
Set-TextValue 'D34' '1.877'
$ws.Range("E34").Value = '  -0.56%  '
Set-TextValue 'D35' '1.168'
Set-TextValue 'D36' '0.7201'
$ws.Range("E36").Value = '  -1.80%  '
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").Value = '1.241.24'
$ws.Range("E38").Value = '  +2.10%  '
Set-TextValue 'D39' '2.797'
$ws.Range("E39").Value = '  -0.72%  '
Set-TextValue 'D40' '0.01802'
$ws.Range("E40").Value = '  +0.53%  '
Set-TextValue 'D41' '0.9090'
$ws.Range("E41").Value = '  -0.94%  '
Set-TextValue 'D42' '6.151'
$ws.Range("E42").Value = '  -3.67%  '
$ws.Range("D43").Value = '2.054.13'
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D45' '67.21'
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D46' '101.08'
$ws.Range("E46").Value = '  -0.92%  '
Set-TextValue 'D47' '7.426'
$ws.Range("E47").Value = '  +10.21%  '
$ws.Range("E48").Value = '  -0.47%  '
Set-TextValue 'D49' '0.4048'
$ws.Range("E49").Value = '  -0.74%  '
Set-TextValue 'D50' '9.165'
$ws.Range("E50").Value = '  -0.49%  '
Set-TextValue 'D51' '1.698'
$ws.Range("E51").Value = '  +2.12%  '
